$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 239
$ws1.Range("F5").Value = 2801
$ws1.Range("F6").Value = 1956
$ws1.Range("F7").Value = 377
$ws1.Range("F8").Value = 125
$ws1.Range("F9").Value = 1015
$ws1.Range("F11").Value = 77
$ws1.Range("F12").Value = 31

# Sheet "全部类型" (All types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 239
$ws4.Range("F5").Value = 2801
$ws4.Range("F6").Value = 1956
$ws4.Range("F7").Value = 377
$ws4.Range("F9").Value = 125
$ws4.Range("F10").Value = 1015
$ws4.Range("F12").Value = 77
$ws4.Range("F13").Value = 31
